$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the win/loss/tie record for each data row (2 through 51)
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 56   # AD
    $ws.Cells.Item($row, 31).Value = 106  # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
